$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "27.163.17"
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("D3").Value = "1.640.37"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "216.95"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  +2.32%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "19.94"
$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").Value = "1.870.50"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "1.639.89"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("E15").Value = "  +2.19%  "

$ws.Range("D16").Value = "66.93"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "27.176.87"
$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("D19").Value = "217.27"
$ws.Range("E19").Value = "  -1.14%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "6.93"
$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("E22").Value = "  +3.37%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "146.53"
$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  -0.67%  "

$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  +1.52%  "

$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").Value = "1.303.89"
$ws.Range("E34").Value = "  +3.84%  "

$ws.Range("E35").Value = "  +0.38%  "

$ws.Range("D36").Value = "2.47"
$ws.Range("E36").Value = "  +1.64%  "

$ws.Range("E37").Value = "  -1.23%  "

$ws.Range("D38").Value = "0.548"
$ws.Range("E38").Value = "  +2.47%  "

$ws.Range("D39").Value = "0.858"
$ws.Range("E39").Value = "  +3.15%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").Value = "0.812"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("E42").Value = "  +5.92%  "

$ws.Range("E43").Value = "  -1.95%  "

$ws.Range("D44").Value = "1.780.46"

$ws.Range("D45").Value = "61.77"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").Value = "91.82"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("E47").Value = "  +1.98%  "

$ws.Range("E48").Value = "  +1.52%  "

$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").Value = "7.64"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("E51").Value = "  +0.12%  "
